$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blood-glucose log continues past the previous last row (625).
# Column B values look numeric, so Excel would otherwise coerce them to
# real numbers on assignment; format the target range as Text first so they
# land as shared strings (matching the source data), then drop the format
# back to Normal so no stray number format lingers on the cells.
$dataRange = $ws.Range("B626:B699")
$dataRange.NumberFormat = "@"

$ws.Range("A626").Value = "2026/02/13 18:13"
$ws.Range("B626").Value = "16.2"
$ws.Range("A627").Value = "2026/02/13 18:18"
$ws.Range("B627").Value = "16.1"
$ws.Range("A628").Value = "2026/02/13 18:23"
$ws.Range("B628").Value = "15.7"
$ws.Range("A629").Value = "2026/02/13 18:28"
$ws.Range("B629").Value = "15.3"
$ws.Range("A630").Value = "2026/02/13 18:33"
$ws.Range("B630").Value = "15.6"
$ws.Range("A631").Value = "2026/02/13 18:38"
$ws.Range("B631").Value = "15.7"
$ws.Range("A632").Value = "2026/02/13 18:43"
$ws.Range("B632").Value = "15.0"
$ws.Range("A633").Value = "2026/02/13 18:48"
$ws.Range("B633").Value = "14.5"
$ws.Range("A634").Value = "2026/02/13 18:53"
$ws.Range("B634").Value = "15.1"
$ws.Range("A635").Value = "2026/02/13 18:58"
$ws.Range("B635").Value = "14.9"
$ws.Range("A636").Value = "2026/02/13 19:03"
$ws.Range("B636").Value = "14.7"
$ws.Range("A637").Value = "2026/02/13 19:08"
$ws.Range("B637").Value = "14.6"
$ws.Range("A638").Value = "2026/02/13 19:13"
$ws.Range("B638").Value = "15.1"
$ws.Range("A639").Value = "2026/02/13 19:18"
$ws.Range("B639").Value = "14.5"
$ws.Range("A640").Value = "2026/02/13 19:23"
$ws.Range("B640").Value = "14.7"
$ws.Range("A641").Value = "2026/02/13 19:28"
$ws.Range("B641").Value = "14.4"
$ws.Range("A642").Value = "2026/02/13 19:33"
$ws.Range("B642").Value = "14.6"
$ws.Range("A643").Value = "2026/02/13 19:38"
$ws.Range("B643").Value = "14.2"
$ws.Range("A644").Value = "2026/02/13 19:43"
$ws.Range("B644").Value = "13.1"
$ws.Range("A645").Value = "2026/02/13 19:48"
$ws.Range("B645").Value = "13.8"
$ws.Range("A646").Value = "2026/02/13 19:53"
$ws.Range("B646").Value = "13.1"
$ws.Range("A647").Value = "2026/02/13 19:58"
$ws.Range("B647").Value = "13.5"
$ws.Range("A648").Value = "2026/02/13 20:03"
$ws.Range("B648").Value = "13.5"
$ws.Range("A649").Value = "2026/02/13 20:08"
$ws.Range("B649").Value = "13.6"
$ws.Range("A650").Value = "2026/02/13 20:13"
$ws.Range("B650").Value = "13.6"
$ws.Range("A651").Value = "2026/02/13 20:18"
$ws.Range("B651").Value = "14.1"
$ws.Range("A652").Value = "2026/02/13 20:23"
$ws.Range("B652").Value = "15.1"
$ws.Range("A653").Value = "2026/02/13 20:28"
$ws.Range("B653").Value = "15.4"
$ws.Range("A654").Value = "2026/02/13 20:33"
$ws.Range("B654").Value = "16.7"
$ws.Range("A655").Value = "2026/02/13 20:38"
$ws.Range("B655").Value = "17.8"
$ws.Range("A656").Value = "2026/02/13 20:43"
$ws.Range("B656").Value = "17.9"
$ws.Range("A657").Value = "2026/02/13 20:48"
$ws.Range("B657").Value = "18.4"
$ws.Range("A658").Value = "2026/02/13 20:53"
$ws.Range("B658").Value = "18.8"
$ws.Range("A659").Value = "2026/02/13 20:58"
$ws.Range("B659").Value = "19.3"
$ws.Range("A660").Value = "2026/02/13 21:03"
$ws.Range("B660").Value = "20.1"
$ws.Range("A661").Value = "2026/02/13 21:08"
$ws.Range("B661").Value = "21.0"
$ws.Range("A662").Value = "2026/02/13 21:13"
$ws.Range("B662").Value = "20.6"
$ws.Range("A663").Value = "2026/02/13 21:18"
$ws.Range("B663").Value = "21.5"
$ws.Range("A664").Value = "2026/02/13 21:23"
$ws.Range("B664").Value = "22.0"
$ws.Range("A665").Value = "2026/02/13 21:28"
$ws.Range("B665").Value = "21.5"
$ws.Range("A666").Value = "2026/02/13 21:33"
$ws.Range("B666").Value = "21.5"
$ws.Range("A667").Value = "2026/02/13 21:38"
$ws.Range("B667").Value = "21.9"
$ws.Range("A668").Value = "2026/02/13 21:43"
$ws.Range("B668").Value = "23.0"
$ws.Range("A669").Value = "2026/02/13 21:48"
$ws.Range("B669").Value = "23.0"
$ws.Range("A670").Value = "2026/02/13 21:53"
$ws.Range("B670").Value = "23.1"
$ws.Range("A671").Value = "2026/02/13 21:58"
$ws.Range("B671").Value = "22.5"
$ws.Range("A672").Value = "2026/02/13 22:03"
$ws.Range("B672").Value = "22.4"
$ws.Range("A673").Value = "2026/02/13 22:08"
$ws.Range("B673").Value = "22.0"
$ws.Range("A674").Value = "2026/02/13 22:13"
$ws.Range("B674").Value = "22.2"
$ws.Range("A675").Value = "2026/02/13 22:18"
$ws.Range("B675").Value = "21.5"
$ws.Range("A676").Value = "2026/02/13 22:23"
$ws.Range("B676").Value = "21.4"
$ws.Range("A677").Value = "2026/02/13 22:28"
$ws.Range("B677").Value = "21.3"
$ws.Range("A678").Value = "2026/02/13 22:33"
$ws.Range("B678").Value = "21.1"
$ws.Range("A679").Value = "2026/02/13 22:38"
$ws.Range("B679").Value = "21.0"
$ws.Range("A680").Value = "2026/02/13 22:43"
$ws.Range("B680").Value = "20.6"
$ws.Range("A681").Value = "2026/02/13 22:48"
$ws.Range("B681").Value = "21.2"
$ws.Range("A682").Value = "2026/02/13 22:53"
$ws.Range("B682").Value = "20.9"
$ws.Range("A683").Value = "2026/02/13 22:58"
$ws.Range("B683").Value = "21.0"
$ws.Range("A684").Value = "2026/02/13 23:03"
$ws.Range("B684").Value = "21.7"
$ws.Range("A685").Value = "2026/02/13 23:08"
$ws.Range("B685").Value = "22.2"
$ws.Range("A686").Value = "2026/02/13 23:13"
$ws.Range("B686").Value = "22.0"
$ws.Range("A687").Value = "2026/02/13 23:18"
$ws.Range("B687").Value = "22.0"
$ws.Range("A688").Value = "2026/02/13 23:23"
$ws.Range("B688").Value = "21.9"
$ws.Range("A689").Value = "2026/02/13 23:28"
$ws.Range("B689").Value = "21.4"
$ws.Range("A690").Value = "2026/02/13 23:33"
$ws.Range("B690").Value = "21.2"
$ws.Range("A691").Value = "2026/02/13 23:38"
$ws.Range("B691").Value = "19.9"
$ws.Range("A692").Value = "2026/02/13 23:43"
$ws.Range("B692").Value = "19.8"
$ws.Range("A693").Value = "2026/02/13 23:48"
$ws.Range("B693").Value = "19.5"
$ws.Range("A694").Value = "2026/02/13 23:53"
$ws.Range("B694").Value = "19.3"
$ws.Range("A695").Value = "2026/02/13 23:58"
$ws.Range("B695").Value = "19.2"
$ws.Range("A696").Value = "2026/02/14 00:03"
$ws.Range("B696").Value = "19.2"
$ws.Range("A697").Value = "2026/02/14 00:08"
$ws.Range("B697").Value = "18.8"
$ws.Range("A698").Value = "2026/02/14 00:13"
$ws.Range("B698").Value = "18.4"
$ws.Range("A699").Value = "2026/02/14 00:18"
$ws.Range("B699").Value = "18.4"

$dataRange.Style = "Normal"

